$wb = $excel.ActiveWorkbook

# --- Sheet "papers": update the accepted-date of an existing entry ---
$papers = $wb.Worksheets.Item("papers")
$papers.Range("K83").Value = 43468

# --- Sheet "papers": append the newly accepted paper (moved from "submitted") ---
$papersTable = $papers.ListObjects.Item(1)
$papersTable.ListRows.Add() | Out-Null

# Copy the formatting of the previous last row onto the freshly added row
$papers.Range("A83:R83").Copy() | Out-Null
$papers.Range("A84:R84").PasteSpecial(-4122) | Out-Null
$papers.Application.CutCopyMode = 0

$papers.Range("A84").Value = "Economic impact of bovine cysticercosis and taeniosis caused by Taenia saginata in Belgium"
$papers.Range("B84").Value = "Jansen, Famke; Dorny, Pierre; Trevisan, Chiara; Dermauw, Veronique; Laranjo-González, Minerva; Allepuz, Alberto; Dupuy, Céline; Krit, Meryam; Gabriël, Sarah; Devleesschauwer, Brecht"
$papers.Range("C84").Value = "Parasites & Vectors"
$papers.Range("D84").Value = "Parasit. Vectors"
$papers.Range("E84").Value = 2018
$papers.Range("F84").Value = "NA"
$papers.Range("G84").Value = "NA"
$papers.Range("H84").Value = "NA"
$papers.Range("I84").Value = "NA"
$papers.Range("J84").Value = "NA"
$papers.Range("K84").Value = 43469
$papers.Range("L84").Value = "A1"
$papers.Range("M84").Value = ""
$papers.Range("N84").Value = ""
$papers.Range("O84").Value = ""
$papers.Range("P84").Value = ""
$papers.Range("Q84").Value = ""
$papers.Range("R84").Value = ""

$papers.Activate()
$papers.Range("A84").Select() | Out-Null

# --- Sheet "submitted": remove the entry that just got accepted ---
$submitted = $wb.Worksheets.Item("submitted")
$submitted.Rows.Item(7).Delete() | Out-Null
